# Fix formatting of scraped floating point numbers and name separators
# (commit: "fix: fixed formatting when scrapping floating point numbers")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$scratch = $ws.Range("Z1")

function Set-TextValue($cellAddr, $text) {
    $formulaText = $text -replace '"', '""'
    $scratch.Formula = '="' + $formulaText + '"'
    $scratch.Copy()
    $ws.Range($cellAddr).PasteSpecial(-4163)
}

# --- Fix provider/proveedor names: replace comma separators with periods ---
Set-TextValue "E23" 'TRABICHET MARIA. VERGARA ADEL Y OTRA'
Set-TextValue "F23" 'TRABICHET MARIA. VERGARA ADEL Y OTRA'
Set-TextValue "E98" 'TRABICHET MARIA. VERGARA ADEL Y OTRA'
Set-TextValue "F98" 'TRABICHET MARIA. VERGARA ADEL Y OTRA'
Set-TextValue "E108" 'TRABICHET MARIA. VERGARA ADEL Y OTRA'
Set-TextValue "F108" 'TRABICHET MARIA. VERGARA ADEL Y OTRA'
Set-TextValue "E113" 'TRABICHET MARIA. VERGARA ADEL Y OTRA'
Set-TextValue "F113" 'TRABICHET MARIA. VERGARA ADEL Y OTRA'
Set-TextValue "E30" 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
Set-TextValue "F30" 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
Set-TextValue "E67" 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
Set-TextValue "F67" 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
Set-TextValue "E80" 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
Set-TextValue "F80" 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
Set-TextValue "E162" 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
Set-TextValue "F162" 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
Set-TextValue "E31" 'MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO'
Set-TextValue "E81" 'MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO'
Set-TextValue "E157" 'MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO'
Set-TextValue "E52" 'SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH'
Set-TextValue "E115" 'SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH'
Set-TextValue "E164" 'SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH'
Set-TextValue "E180" 'SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH'
Set-TextValue "E120" 'RICCOTTI. MARIANA EDITH'
Set-TextValue "E166" 'ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN'
Set-TextValue "F126" 'MERCANZINI. GASTON ARIEL'

# --- Fix Importe (amount) values: convert "1.234,56" (es-AR) style text to "1234.56" ---
Set-TextValue "H2" '912.00'
Set-TextValue "H3" '3500.00'
Set-TextValue "H4" '3220.00'
Set-TextValue "H5" '96000.00'
Set-TextValue "H6" '78000.00'
Set-TextValue "H7" '30625.00'
Set-TextValue "H8" '95000.00'
Set-TextValue "H9" '160930.00'
Set-TextValue "H10" '814.56'
Set-TextValue "H11" '437.58'
Set-TextValue "H12" '4199.92'
Set-TextValue "H13" '355276.64'
Set-TextValue "H14" '43184.57'
Set-TextValue "H15" '6277.50'
Set-TextValue "H16" '5515.00'
Set-TextValue "H17" '2207.09'
Set-TextValue "H18" '14064.02'
Set-TextValue "H19" '10202.40'
Set-TextValue "H20" '4200.00'
Set-TextValue "H21" '7460.00'
Set-TextValue "H22" '1400.00'
Set-TextValue "H23" '114.00'
Set-TextValue "H24" '2092.50'
Set-TextValue "H25" '265.35'
Set-TextValue "H26" '11000.00'
Set-TextValue "H27" '7800.00'
Set-TextValue "H28" '8030.90'
Set-TextValue "H29" '708.94'
Set-TextValue "H30" '1015.44'
Set-TextValue "H31" '705.00'
Set-TextValue "H32" '1268.00'
Set-TextValue "H33" '13311.00'
Set-TextValue "H34" '306.00'
Set-TextValue "H35" '590.00'
Set-TextValue "H36" '105000.00'
Set-TextValue "H37" '81462.65'
Set-TextValue "H38" '1007.70'
Set-TextValue "H39" '1850.00'
Set-TextValue "H40" '50.00'
Set-TextValue "H41" '45.14'
Set-TextValue "H42" '3642.00'
Set-TextValue "H43" '9760.00'
Set-TextValue "H44" '14386.90'
Set-TextValue "H45" '3072.71'
Set-TextValue "H46" '907.50'
Set-TextValue "H47" '1460.00'
Set-TextValue "H48" '4.96'
Set-TextValue "H49" '10520.66'
Set-TextValue "H50" '6646.92'
Set-TextValue "H51" '1132.00'
Set-TextValue "H52" '1280.00'
Set-TextValue "H53" '22585.96'
Set-TextValue "H54" '360.19'
Set-TextValue "H55" '101065.22'
Set-TextValue "H56" '864.63'
Set-TextValue "H57" '155.00'
Set-TextValue "H58" '6537.00'
Set-TextValue "H59" '5964.75'
Set-TextValue "H60" '640.00'
Set-TextValue "H61" '43672.00'
Set-TextValue "H62" '3153.00'
Set-TextValue "H63" '305.00'
Set-TextValue "H64" '1400.00'
Set-TextValue "H65" '120.00'
Set-TextValue "H66" '12860.00'
Set-TextValue "H67" '5.04'
Set-TextValue "H68" '840.00'
Set-TextValue "H69" '38889.51'
Set-TextValue "H70" '400.00'
Set-TextValue "H71" '388.00'
Set-TextValue "H72" '19056.00'
Set-TextValue "H73" '760.00'
Set-TextValue "H74" '178.38'
Set-TextValue "H75" '396.00'
Set-TextValue "H76" '1376.85'
Set-TextValue "H77" '489.00'
Set-TextValue "H78" '26495.00'
Set-TextValue "H79" '1337.05'
Set-TextValue "H80" '109.20'
Set-TextValue "H81" '2871.00'
Set-TextValue "H82" '8159.00'
Set-TextValue "H83" '157.58'
Set-TextValue "H84" '1570.00'
Set-TextValue "H85" '118.46'
Set-TextValue "H86" '44400.00'
Set-TextValue "H87" '1050.00'
Set-TextValue "H88" '0.64'
Set-TextValue "H89" '21221.88'
Set-TextValue "H90" '12000.00'
Set-TextValue "H91" '410745.87'
Set-TextValue "H92" '308.60'
Set-TextValue "H93" '289935.82'
Set-TextValue "H94" '55047.50'
Set-TextValue "H95" '1.24'
Set-TextValue "H96" '2100.00'
Set-TextValue "H97" '8.94'
Set-TextValue "H98" '2520.00'
Set-TextValue "H99" '177.64'
Set-TextValue "H100" '1825.00'
Set-TextValue "H101" '14.24'
Set-TextValue "H102" '4400.00'
Set-TextValue "H103" '2868.72'
Set-TextValue "H104" '3789.00'
Set-TextValue "H105" '21.38'
Set-TextValue "H106" '15085.00'
Set-TextValue "H107" '468.00'
Set-TextValue "H108" '385.00'
Set-TextValue "H109" '1760.00'
Set-TextValue "H110" '2400.00'
Set-TextValue "H111" '9429.10'
Set-TextValue "H112" '1220.21'
Set-TextValue "H113" '36.00'
Set-TextValue "H114" '664.00'
Set-TextValue "H115" '800.00'
Set-TextValue "H116" '90.00'
Set-TextValue "H117" '450.00'
Set-TextValue "H118" '2310.00'
Set-TextValue "H119" '650.00'
Set-TextValue "H120" '1000.00'
Set-TextValue "H121" '540.00'
Set-TextValue "H122" '590.00'
Set-TextValue "H123" '5000.00'
Set-TextValue "H124" '2000.00'
Set-TextValue "H125" '1485.00'
Set-TextValue "H126" '6000.00'
Set-TextValue "H127" '4200.00'
Set-TextValue "H128" '1250.00'
Set-TextValue "H129" '102.80'
Set-TextValue "H130" '840.50'
Set-TextValue "H131" '1712.78'
Set-TextValue "H132" '77.00'
Set-TextValue "H133" '1476.65'
Set-TextValue "H134" '154000.00'
Set-TextValue "H135" '10234.88'
Set-TextValue "H136" '2300.00'
Set-TextValue "H137" '1450.00'
Set-TextValue "H138" '1800.00'
Set-TextValue "H139" '1000.00'
Set-TextValue "H140" '10167.50'
Set-TextValue "H141" '4000.00'
Set-TextValue "H142" '700.00'
Set-TextValue "H143" '800.00'
Set-TextValue "H144" '1000.00'
Set-TextValue "H145" '17608.51'
Set-TextValue "H146" '2500.00'
Set-TextValue "H147" '950.00'
Set-TextValue "H148" '1000.00'
Set-TextValue "H149" '6030.66'
Set-TextValue "H150" '2000.00'
Set-TextValue "H151" '600.00'
Set-TextValue "H152" '320.00'
Set-TextValue "H153" '15940.00'
Set-TextValue "H154" '1300.00'
Set-TextValue "H155" '301.29'
Set-TextValue "H156" '1430.00'
Set-TextValue "H157" '220.00'
Set-TextValue "H158" '565.00'
Set-TextValue "H159" '180.00'
Set-TextValue "H160" '150.00'
Set-TextValue "H161" '4402.00'
Set-TextValue "H162" '691.71'
Set-TextValue "H163" '5160.00'
Set-TextValue "H164" '740.00'
Set-TextValue "H165" '277.04'
Set-TextValue "H166" '2705.00'
Set-TextValue "H167" '858.90'
Set-TextValue "H168" '36.00'
Set-TextValue "H169" '15061.00'
Set-TextValue "H170" '1600.00'
Set-TextValue "H171" '1406.00'
Set-TextValue "H172" '27781.60'
Set-TextValue "H173" '80.25'
Set-TextValue "H174" '6.52'
Set-TextValue "H175" '1400.00'
Set-TextValue "H176" '1082.56'
Set-TextValue "H177" '991.00'
Set-TextValue "H178" '4223.40'
Set-TextValue "H179" '843.00'
Set-TextValue "H180" '10207.00'
Set-TextValue "H181" '79.20'
Set-TextValue "H182" '3578.00'
Set-TextValue "H183" '1209.60'
Set-TextValue "H184" '2721.00'
Set-TextValue "H185" '130.00'
Set-TextValue "H186" '420.00'
Set-TextValue "H187" '22.50'
Set-TextValue "H188" '87.00'
Set-TextValue "H189" '17600.00'
Set-TextValue "H190" '219.00'
Set-TextValue "H191" '3240.00'
Set-TextValue "H192" '62658.00'
Set-TextValue "H193" '6440.00'
Set-TextValue "H194" '7.50'
Set-TextValue "H195" '1361.16'
Set-TextValue "H196" '112.00'
Set-TextValue "H197" '615.00'
Set-TextValue "H198" '274069.66'
Set-TextValue "H199" '7036.39'
Set-TextValue "H200" '4977.17'
Set-TextValue "H201" '728609.27'
Set-TextValue "H202" '2150.00'
Set-TextValue "H203" '3635.67'
Set-TextValue "H204" '750.00'
Set-TextValue "H205" '480.00'
Set-TextValue "H206" '1200.00'
Set-TextValue "H207" '291000.00'
Set-TextValue "H208" '410385.00'
Set-TextValue "H209" '60000.00'
Set-TextValue "H210" '257200.00'
Set-TextValue "H211" '115000.00'
Set-TextValue "H212" '206266.00'
Set-TextValue "H213" '446000.00'
Set-TextValue "H214" '223000.00'
Set-TextValue "H215" '223000.00'
Set-TextValue "H216" '160752.00'
Set-TextValue "H217" '121000.00'
Set-TextValue "H218" '290.00'
Set-TextValue "H219" '8380.00'
Set-TextValue "H220" '54000.00'
Set-TextValue "H221" '2544.00'
Set-TextValue "H222" '325353.00'
Set-TextValue "H223" '700.00'
Set-TextValue "H224" '400.00'
Set-TextValue "H225" '890.00'
Set-TextValue "H226" '3750.00'

$scratch.Clear()
Write-Output "done"
